$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '42.315.33'
$ws.Range('E2').Value = '  +1.31%  '
Set-TextValue $ws 'D3' '2.272.40'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue $ws 'D5' '306.86'
$ws.Range('E5').Value = '  +0.95%  '
Set-TextValue $ws 'D6' '97.29'
$ws.Range('E6').Value = '  +4.97%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +1.36%  '
Set-TextValue $ws 'D10' '35.42'
$ws.Range('E10').Value = '  +8.64%  '
Set-TextValue $ws 'D11' '0.0795'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('E13').Value = '  -0.29%  '
Set-TextValue $ws 'D14' '2.621.57'
$ws.Range('E14').Value = '  -0.08%  '
Set-TextValue $ws 'D15' '14.36'
$ws.Range('E15').Value = '  +0.34%  '
Set-TextValue $ws 'D16' '2.272.18'
$ws.Range('E16').Value = '  -0.22%  '
Set-TextValue $ws 'D17' '0.794'
$ws.Range('E17').Value = '  +2.21%  '
Set-TextValue $ws 'D18' '42.205.01'
$ws.Range('E18').Value = '  +1.24%  '
Set-TextValue $ws 'D19' '12.48'
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('E21').Value = '  +0.13%  '
Set-TextValue $ws 'D22' '67.53'
$ws.Range('E22').Value = '  +0.55%  '
Set-TextValue $ws 'D23' '240.25'
$ws.Range('E23').Value = '  -1.38%  '
$ws.Range('E24').Value = '  +0.42%  '
Set-TextValue $ws 'D25' '1.94'
$ws.Range('E25').Value = '  +0.82%  '
Set-TextValue $ws 'D26' '0.999'
$ws.Range('E26').Value = '  -0.15%  '
Set-TextValue $ws 'D27' '23.81'
$ws.Range('E27').Value = '  -0.67%  '
Set-TextValue $ws 'D28' '37.42'
$ws.Range('E28').Value = '  +5.88%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +1.61%  '
Set-TextValue $ws 'D31' '159.14'
$ws.Range('E31').Value = '  -0.85%  '
Set-TextValue $ws 'D32' '5.25'
$ws.Range('E32').Value = '  +0.31%  '
Set-TextValue $ws 'D33' '0.999'
$ws.Range('E33').Value = '  -0.04%  '
Set-TextValue $ws 'D34' '3.15'
$ws.Range('E34').Value = '  +4.54%  '
Set-TextValue $ws 'D35' '0.0740'
$ws.Range('E35').Value = '  -0.48%  '
Set-TextValue $ws 'D36' '17.00'
$ws.Range('E36').Value = '  +0.60%  '
Set-TextValue $ws 'D37' '2.37'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('E39').Value = '  +1.83%  '
$ws.Range('E40').Value = '  -1.58%  '
Set-TextValue $ws 'D41' '4.07'
$ws.Range('E41').Value = '  +4.03%  '
$ws.Range('E42').Value = '  +13.39%  '
Set-TextValue $ws 'D43' '1.996.33'
$ws.Range('E43').Value = '  -0.57%  '
$ws.Range('E44').Value = '  +0.89%  '
Set-TextValue $ws 'D45' '18.77'
$ws.Range('E45').Value = '  -3.73%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D46' '2.94'
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D47' '9.96'
$ws.Range('E47').Value = '  -3.34%  '
Set-TextValue $ws 'D48' '52.86'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('E49').Value = '  +0.52%  '
Set-TextValue $ws 'D50' '72.08'
$ws.Range('E50').Value = '  -0.09%  '
Set-TextValue $ws 'D51' '91.48'
$ws.Range('E51').Value = '  +0.42%  '
